$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (ligand/receptor cluster text stays the same,
# but numeric statistics are recalculated with the new TPM-based figures)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03315566666666667
$ws.Range("H2").Value = 0.099467
$ws.Range("I2").Value = 0.1557603470145164
$ws.Range("J2").Value = 0.1557603470145164
$ws.Range("Q2").Value = 0.3818614830108889
$ws.Range("R2").Value = 3.436753347098
$ws.Range("S2").Value = 0.1557603470145164
$ws.Range("T2").Value = 0.1557603470145164

# Add new row 3 for sending cluster "MuSCs"
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl3"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1797076666666667
$ws.Range("H3").Value = 0.539123
$ws.Range("I3").Value = 0.8442396529854836
$ws.Range("J3").Value = 0.8442396529854836
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.51723133333333
$ws.Range("N3").Value = 34.551694
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 2.069734769373555
$ws.Range("R3").Value = 18.627612924362
$ws.Range("S3").Value = 0.8442396529854836
$ws.Range("T3").Value = 0.8442396529854836
